# Apply retriever evaluation results update:
# - Refresh accuracy/MRR metrics for sentence-transformers/all-MiniLM-L6-v2 (row 2)
# - Refresh accuracy/MRR metrics for dunzhang/stella_en_1.5B_v5 load_in_8bit=True (row 5)
# - Replace dunzhang/stella_en_1.5B_v5 (load_in_8bit=False) rows 6 & 11 with
#   amazon.titan-embed-text-v2:0 results (AWS Bedrock Titan embedding model)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: sentence-transformers/all-MiniLM-L6-v2 (HuggingFace QA Dataset) ---
$ws.Range("D2").Value = 0.6307692307692307
$ws.Range("E2").Value = 0.7076923076923077
$ws.Range("F2").Value = 0.7230769230769231
$ws.Range("G2").Value = 0.7384615384615385
$ws.Range("H2").Value = 0.7846153846153846
$ws.Range("I2").Value = 0.6874358974358974

# --- Row 5: dunzhang/stella_en_1.5B_v5, load_in_8bit=True (HuggingFace QA Dataset) ---
$ws.Range("D5").Value = 0.2923076923076923
$ws.Range("E5").Value = 0.6153846153846154
$ws.Range("F5").Value = 0.7076923076923077
$ws.Range("G5").Value = 0.7384615384615385
$ws.Range("H5").Value = 0.7538461538461538
$ws.Range("I5").Value = 0.4953846153846154

# --- Row 6: replace stella_en_1.5B_v5 (load_in_8bit=False) with amazon.titan-embed-text-v2:0 (HuggingFace QA Dataset) ---
$ws.Range("A6").Value = "amazon.titan-embed-text-v2:0"
$ws.Range("C6").Value = "{'name': 'amazon.titan-embed-text-v2:0', 'instruction': 'Instruct: Represent this passage for retrieval in response to relevant questions.`nQuery:', 'query_instruction': 'Instruct: Given a query, find the most relevant passages that can provide the answer.`nPassage:', 'model_kwargs': {'aws': True, 'aws_creds_file': '/home/ubuntu/Multi-Agent-LLM-System-with-LangGraph-RAG-and-LangChain/config/config.ini', 'aws_config_name': 'BedRock_LLM_API'}}"
$ws.Range("D6").Value = 0.8153846153846154
$ws.Range("E6").Value = 0.9076923076923077
$ws.Range("F6").Value = 0.9538461538461539
$ws.Range("G6").Value = 0.9538461538461539
$ws.Range("H6").Value = 0.9692307692307692
$ws.Range("I6").Value = 0.88

# --- Row 11: replace stella_en_1.5B_v5 (load_in_8bit=False) with amazon.titan-embed-text-v2:0 (PubMed filtered Dataset) ---
$ws.Range("A11").Value = "amazon.titan-embed-text-v2:0"
$ws.Range("C11").Value = "{'name': 'amazon.titan-embed-text-v2:0', 'instruction': 'Instruct: Represent this passage for retrieval in response to relevant questions.`nQuery:', 'query_instruction': 'Instruct: Given a query, find the most relevant passages that can provide the answer.`nPassage:', 'model_kwargs': {'aws': True, 'aws_creds_file': '/home/ubuntu/Multi-Agent-LLM-System-with-LangGraph-RAG-and-LangChain/config/config.ini', 'aws_config_name': 'BedRock_LLM_API'}}"
$ws.Range("D11").Value = 0.9230769230769231
$ws.Range("E11").Value = 0.9230769230769231
$ws.Range("F11").Value = 0.9230769230769231
$ws.Range("G11").Value = 0.9230769230769231
$ws.Range("H11").Value = 0.9230769230769231
$ws.Range("I11").Value = 0.9230769230769231

$wb.Save()
